$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

# Row 2
Set-TextCell "D2" "259.38"
Set-TextCell "E2" "0.96%"
Set-TextCell "G2" "23"

# Row 3
Set-TextCell "D3" "26.87"
Set-TextCell "E3" "-0.16%"
Set-TextCell "G3" "23"

# Row 4
Set-TextCell "D4" "4.672"
Set-TextCell "E4" "0.79%"
Set-TextCell "G4" "23"

# Row 5
Set-TextCell "D5" "0.06059"
Set-TextCell "E5" "3.09%"
Set-TextCell "G5" "23"

# Row 6
Set-TextCell "D6" "6.689"
Set-TextCell "E6" "0.85%"
Set-TextCell "G6" "23"

# Row 7
Set-TextCell "D7" "0.8616"
Set-TextCell "E7" "0.42%"
Set-TextCell "G7" "23"

# Row 8
Set-TextCell "D8" "0.9220"
Set-TextCell "E8" "-1.88%"
Set-TextCell "G8" "23"

# Row 9
Set-TextCell "E9" "-0.81%"
Set-TextCell "G9" "23"

# Row 10
Set-TextCell "D10" "0.05203"
Set-TextCell "E10" "24.48%"
Set-TextCell "G10" "23"

# Row 11
Set-TextCell "D11" "0.07094"
Set-TextCell "E11" "0.00%"
Set-TextCell "G11" "23"

# Row 12
Set-TextCell "D12" "0.03123"
Set-TextCell "E12" "-0.86%"
Set-TextCell "G12" "23"

# Row 13
Set-TextCell "D13" "0.09130"
Set-TextCell "E13" "-0.30%"
Set-TextCell "G13" "23"

# Row 14
Set-TextCell "D14" "0.001545"
Set-TextCell "E14" "0.42%"
Set-TextCell "G14" "23"

# Row 15
Set-TextCell "D15" "0.0006064"
Set-TextCell "E15" "-94.19%"
Set-TextCell "G15" "23"

# Row 16
Set-TextCell "D16" "0.006005"
Set-TextCell "E16" "-3.57%"
Set-TextCell "G16" "23"

# Row 17
Set-TextCell "D17" "3.484"
Set-TextCell "E17" "-1.01%"
Set-TextCell "G17" "23"

# Row 18
Set-TextCell "D18" "3.169"
Set-TextCell "E18" "-1.18%"
Set-TextCell "G18" "23"

# Row 19
Set-TextCell "E19" "-1.30%"
Set-TextCell "G19" "23"

# Row 20
Set-TextCell "E20" "2.46%"
Set-TextCell "G20" "23"

# Row 21
Set-TextCell "E21" "-0.16%"
Set-TextCell "G21" "23"

# Row 22
Set-TextCell "D22" "4.086"
Set-TextCell "E22" "6.48%"
Set-TextCell "G22" "23"

# Row 23
Set-TextCell "D23" "0.04238"
Set-TextCell "E23" "0.19%"
Set-TextCell "G23" "23"

# Row 24
Set-TextCell "D24" "0.001216"
Set-TextCell "E24" "-0.72%"
Set-TextCell "G24" "23"

# Row 25
Set-TextCell "D25" "0.004021"
Set-TextCell "G25" "23"

# Row 26
Set-TextCell "E26" "-0.08%"
Set-TextCell "G26" "23"

# Row 27
Set-TextCell "G27" "23"

# Row 28
Set-TextCell "G28" "23"

# Row 29
Set-TextCell "G29" "23"

# Row 30
Set-TextCell "G30" "23"

# Row 31
Set-TextCell "G31" "23"

# Row 32
Set-TextCell "G32" "23"

# Row 33
Set-TextCell "G33" "23"

# Row 34
Set-TextCell "G34" "23"

# Row 35
Set-TextCell "G35" "23"

# Row 36
Set-TextCell "G36" "23"

# Row 37
Set-TextCell "G37" "23"

# Row 38
Set-TextCell "G38" "23"

# Row 39
Set-TextCell "G39" "23"

# Row 40
Set-TextCell "D40" "0.03869"
Set-TextCell "E40" "1.06%"
Set-TextCell "G40" "23"

# Row 41
Set-TextCell "E41" "1.47%"
Set-TextCell "G41" "23"

# Row 42
Set-TextCell "D42" "0.004067"
Set-TextCell "E42" "-34.95%"
Set-TextCell "G42" "23"

# Row 43
Set-TextCell "D43" "0.01492"
Set-TextCell "E43" "30.52%"
Set-TextCell "G43" "23"

# Row 44
Set-TextCell "E44" "-0.08%"
Set-TextCell "G44" "23"

# Row 45
Set-TextCell "D45" "0.00005187"
Set-TextCell "E45" "-5.01%"
Set-TextCell "G45" "23"

# Row 46
Set-TextCell "E46" "-0.08%"
Set-TextCell "G46" "23"

# Row 47
Set-TextCell "B47" "CoinbaseStockToken"
Set-TextCell "C47" "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextCell "D47" "0.05453"
Set-TextCell "E47" "-27.13%"
Set-TextCell "G47" "23"

# Row 48
Set-TextCell "B48" "BOLO"
Set-TextCell "C48" "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextCell "D48" "0.1353"
Set-TextCell "E48" "-42.84%"
Set-TextCell "G48" "23"

# Row 49
Set-TextCell "E49" "-0.08%"
Set-TextCell "G49" "23"

# Row 50
Set-TextCell "E50" "-0.08%"
Set-TextCell "G50" "23"

# Row 51
Set-TextCell "G51" "23"
